$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert two new rows before row 10 (current blank separator row) and
# copy the formatting from row 9 (the last "header block" row) into them.
$ws.Rows("10:11").Insert()
$ws.Range("A9:C9").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122) | Out-Null
$ws.Range("A10:C10").PasteSpecial(-4123) | Out-Null
$ws.Range("A9:C9").Copy()
$ws.Range("A11:C11").PasteSpecial(-4122) | Out-Null
$ws.Range("A11:C11").PasteSpecial(-4123) | Out-Null
$excel.CutCopyMode = $false
$ws.Rows("10:11").RowHeight = 15
$ws.Range("C10:C11").ClearContents()

$ws.Range("A10").Value = "RutaTolerancia"
$ws.Range("A11").Value = "RutaTemplateTol"
$ws.Range("B10").Value = "C:\Users\ROBTIRELEO\Documents\JACK3\lib\net45\Data\Tolerancia"
$ws.Range("B11").Value = "C:\Users\ROBTIRELEO\Documents\JACK3\lib\net45\Data\Templates\Template Tolerancia.xlsx"

# New row appended at the bottom of the settings list (row 45 after the
# insertion above shifted the old list down by two rows). A45 reuses the
# bold/highlighted look of the "Grupo..." cells (now at row 23 after the
# insert), B45 reuses the plain highlighted look of the row above it.
$ws.Range("C23").Copy()
$ws.Range("A45").PasteSpecial(-4122) | Out-Null
$ws.Range("A45").PasteSpecial(-4123) | Out-Null
$ws.Range("B44").Copy()
$ws.Range("B45").PasteSpecial(-4122) | Out-Null
$ws.Range("B45").PasteSpecial(-4123) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A45").Value = "ValorRango"
$ws.Range("B45").Value = 10000

$window = $excel.ActiveWindow
$window.ScrollRow = 1
$ws.Range("A45").Select()
